$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new data rows (id, name, avgmark, accumulationcredit, birthday-serial).
# "S" is a brand-new value for column B (not previously present as a shared string).
$rows = @(
    @(51503002, "S", 5, 55, 43159),
    @(51503003, "S", 5, 5,  43159),
    @(51503004, "S", 5, 5,  43159)
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Give the new date cell the same number format as the existing BirthDay
    # column by copying formats from the row above, so it shares the cellXf
    # (style index) instead of minting a new one.
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = $row[4]
    $ws.Cells.Item($r - 1, 5).Copy() | Out-Null
    $eCell.PasteSpecial(-4122) | Out-Null

    $r = $r + 1
}
$excel.CutCopyMode = $false

# Column A gets an explicit bestFit-style width, like the real edit did.
$ws.Columns.Item(1).ColumnWidth = 9.166666666666666

# Move the active selection to the new last cell, matching the author's edit.
$ws.Range("E6").Select() | Out-Null
